# Fix the 2050 column header label (previously a stray numeric value) and
# remove the "Total" summary row from each scenario table.

$wb = $excel.ActiveWorkbook

function Set-TextLabel {
    param($ws, [string]$cellAddress, [string]$text, [string]$formatSourceAddress)
    $cell = $ws.Range($cellAddress)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    # Restore the original visual style (border/font/alignment) that got
    # reset when the cell was forced to a text number format.
    $ws.Range($formatSourceAddress).Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

# --- Sheets with the A1:E13 layout (label row + 11 data rows + Total row) ---
$ws = $wb.Worksheets.Item("Potencia Acumulada - SIN (MW)")
Set-TextLabel $ws "E1" "2050" "D1"
$ws.Rows.Item(13).Delete()

$ws = $wb.Worksheets.Item("Geracao Periodo Medio (MWMed)")
Set-TextLabel $ws "E1" "2050" "D1"
$ws.Rows.Item(13).Delete()

$ws = $wb.Worksheets.Item("Atendimento a Ponta(MW)")
Set-TextLabel $ws "E1" "2050" "D1"
$ws.Rows.Item(13).Delete()

$ws = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
Set-TextLabel $ws "E1" "2041-2050" "D1"
$ws.Rows.Item(13).Delete()

# --- "Emissoes Totais" sheet: only the label needs fixing, no Total row ---
$ws = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
Set-TextLabel $ws "E1" "2050" "D1"

# --- "Custo Total" sheet: just drop the Total row (no column E here) ---
$ws = $wb.Worksheets.Item('Custo Total (bilhões de R$)')
$ws.Rows.Item(4).Delete()
